$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$tr.Paragraphs(2, 1).Text = ""
$tr.Paragraphs(2, 1).Text = "To align and enrich the data from the structured JSON and unstructured PDF content, let's proceed with the task step-by-step:"

$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "### Step 1: Align Assets"

$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "From the JSON data, we have the asset `"DC1 Ingram Micro`" located at `"DC1, Drayton Way, Apex Park NN11 8NF, Daventry, United Kingdom`". We need to find a corresponding description in the PDF content. However, the provided PDF content does not mention this specific asset or address. Therefore, no direct alignment can be made based on the provided PDF content."

$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "### Step 2: Extract & Enrich"

$tr.Paragraphs(6, 1).Text = ""
$tr.Paragraphs(6, 1).Text = "Since we couldn't find a direct match in the PDF content for `"DC1 Ingram Micro`", we cannot extract specific details like market highlights, investment rationale, risk factors, financials, or physical specs for this asset from the PDF. However, we can summarize the available data from the JSON:"

$tr.Paragraphs(7, 1).Text = ""
$tr.Paragraphs(7, 1).Text = "#### JSON Data Summary for `"DC1 Ingram Micro`":"

$tr.Paragraphs(8, 1).Text = ""
$tr.Paragraphs(8, 1).Text = "- **Asset Name**: DC1 Ingram Micro"

$tr.Paragraphs(9, 1).Text = ""
$tr.Paragraphs(9, 1).Text = "- **Asset Type**: Logistics"

$tr.Paragraphs(10, 1).Text = ""
$tr.Paragraphs(10, 1).Text = "- **Tenure**: Freehold"

$tr.Paragraphs(11, 1).Text = ""
$tr.Paragraphs(11, 1).Text = "- **Address**: DC1, Drayton"
